$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the probe measurements for the fan shroud calibration pass.
$ws.Range("B7").Value = 2124.8000000000002
$ws.Range("D8").Value = 101
$ws.Range("D9").Value = 101
$ws.Range("D10").Value = 101

# Move the active selection to D9, matching the author's last-saved view.
$ws.Range("D9").Select()
